$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Give every "straight connector arrow" shape on slides 1-3 an explicit
#    2pt (25400 EMU) line weight (<a:ln> -> <a:ln w="25400">).
# ---------------------------------------------------------------------
for ($si = 1; $si -le 3; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -like "*Conector recto de flecha*") {
            $sh.Line.Weight = 2
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 3: merge the "Ethereum" + "/RSK" runs into a single run so the
#    leftover spell-check "err" flag from the first run disappears.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "Ethereum/RSK") {
        $tr = $sh.TextFrame.TextRange
        # Drop the first run ("Ethereum") entirely...
        $firstRun = $tr.Characters(1, 8)
        $firstRun.Text = ""
        # ...then re-prepend it onto what remains ("/RSK"), which keeps
        # the second run's (error-free) character formatting.
        $rest = $tr.Characters(1, 4)
        $rest.Text = "Ethereum" + $rest.Text
    }
}

# ---------------------------------------------------------------------
# 3) Slide 9: remove the standalone "Mining" title shape.
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
for ($i = $s9.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s9.Shapes.Item($i)
    if ($sh.Name -eq "Title 3") {
        $sh.Cut()
    }
}
